# Weekly update: insert a new Primera/Segunda price pair for Brócoli
# (Terminal La Palmera de La Serena) right above the existing block,
# shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 491-492, pushing old row 491.. down to 493..
$ws.Range("A491:A492").EntireRow.Insert()

# New row 491 - "Primera" quality
$ws.Cells.Item(491, 1).Value = 8
$ws.Cells.Item(491, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(491, 3).Value = "Coquimbo"
$ws.Cells.Item(491, 4).Value = 44578
$ws.Cells.Item(491, 5).Value = 4
$ws.Cells.Item(491, 6).Value = 100112023
$ws.Cells.Item(491, 7).Value = "Brócoli"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 2560
$ws.Cells.Item(491, 11).Value = 650
$ws.Cells.Item(491, 12).Value = 700
$ws.Cells.Item(491, 13).Value = 675
$ws.Cells.Item(491, 14).Value = "$/unidad"
$ws.Cells.Item(491, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(491, 16).Value = 675
$ws.Cells.Item(491, 17).Value = 1
$ws.Cells.Item(491, 18).Value = "Hortaliza"

# New row 492 - "Segunda" quality
$ws.Cells.Item(492, 1).Value = 8
$ws.Cells.Item(492, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(492, 3).Value = "Coquimbo"
$ws.Cells.Item(492, 4).Value = 44578
$ws.Cells.Item(492, 5).Value = 4
$ws.Cells.Item(492, 6).Value = 100112023
$ws.Cells.Item(492, 7).Value = "Brócoli"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Segunda"
$ws.Cells.Item(492, 10).Value = 1500
$ws.Cells.Item(492, 11).Value = 550
$ws.Cells.Item(492, 12).Value = 600
$ws.Cells.Item(492, 13).Value = 575
$ws.Cells.Item(492, 14).Value = "$/unidad"
$ws.Cells.Item(492, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(492, 16).Value = 575
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"
